# Macroferia Regional de Talca - Uva
# Insert a new weekly record at row 436 (shifting existing rows 436-523 down
# to 437-524), matching the "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 436; this shifts rows 436:523 down to 437:524
# and the new row inherits the formatting (incl. the date number format)
# from the row above it, just like Excel's native Insert behaviour.
$ws.Rows.Item(436).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(436, 1).Value = 5
$ws.Cells.Item(436, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(436, 3).Value = "Maule"
$ws.Cells.Item(436, 4).Value = 44785
$ws.Cells.Item(436, 5).Value = 7
$ws.Cells.Item(436, 6).Value = "Fruta"
$ws.Cells.Item(436, 7).Value = 100109
$ws.Cells.Item(436, 8).Value = "Uva"
$ws.Cells.Item(436, 9).Value = 100109001
$ws.Cells.Item(436, 10).Value = "Uva"
$ws.Cells.Item(436, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(436, 12).Value = "Segunda"
$ws.Cells.Item(436, 13).Value = 300
$ws.Cells.Item(436, 14).Value = 8000
$ws.Cells.Item(436, 15).Value = 8000
$ws.Cells.Item(436, 16).Value = 8000
$ws.Cells.Item(436, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(436, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(436, 19).Value = 1000
$ws.Cells.Item(436, 20).Value = 8
